$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns H, I, AF for rows 2-8 with the new encoded values.
for ($row = 2; $row -le 8; $row++) {
    $ws.Range("H$row").Value = 2147482380
    $ws.Range("I$row").Value = 2147482384
    $ws.Range("AF$row").Value = 2147482376
}
